# "error solve ifrs list"
# The 2014-2018 (IFRS-consolidated) actuals in rows 2-6 were populated with
# the wrong figures; replace them with the corrected SK Securities values,
# including dropping the discontinued FCF (U) column. The 2019E-2021E
# forecast rows (7-9) are removed entirely, leaving only their A-C labels.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: 2014/12  (IFRS연결)
$ws.Range("D2").Value = 5269
$ws.Range("E2").Value = 95
$ws.Range("F2").Value = 95
$ws.Range("G2").Value = 67
$ws.Range("H2").Value = 34
$ws.Range("I2").Value = 34
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 36757
$ws.Range("L2").Value = 32909
$ws.Range("M2").Value = 3847
$ws.Range("N2").Value = 3847
$ws.Range("O2").Value = 0
$ws.Range("P2").Value = 1620
$ws.Range("Q2").Value = 3065
$ws.Range("R2").Value = -46
$ws.Range("S2").Value = -2639
$ws.Range("T2").Value = 8
$ws.Range("V2").Value = 5794
$ws.Range("W2").Value = 1.81
$ws.Range("X2").Value = 0.65
$ws.Range("Y2").Value = 0.9
$ws.Range("Z2").Value = 0.09
$ws.Range("AA2").Value = 855.39
$ws.Range("AB2").Value = 137.47
$ws.Range("AC2").Value = 10
$ws.Range("AD2").Value = 83.28
$ws.Range("AE2").Value = 1142
$ws.Range("AF2").Value = 0.73
$ws.Range("AG2").Value = 0
$ws.Range("AH2").Value = 0
$ws.Range("AI2").Value = 0
$ws.Range("AJ2").Value = 341155134
$ws.Range("U2").ClearContents()

# Row 3: 2015/12  (IFRS연결)
$ws.Range("D3").Value = 4709
$ws.Range("E3").Value = 202
$ws.Range("F3").Value = 202
$ws.Range("G3").Value = 266
$ws.Range("H3").Value = 230
$ws.Range("I3").Value = 230
$ws.Range("K3").Value = 35081
$ws.Range("L3").Value = 30966
$ws.Range("M3").Value = 4114
$ws.Range("N3").Value = 4114
$ws.Range("P3").Value = 1620
$ws.Range("Q3").Value = -448
$ws.Range("R3").Value = 721
$ws.Range("S3").Value = -351
$ws.Range("T3").Value = 13
$ws.Range("V3").Value = 5190
$ws.Range("W3").Value = 4.3
$ws.Range("X3").Value = 4.89
$ws.Range("Y3").Value = 5.78
$ws.Range("Z3").Value = 0.64
$ws.Range("AA3").Value = 752.61
$ws.Range("AB3").Value = 153.96
$ws.Range("AC3").Value = 67
$ws.Range("AD3").Value = 16.16
$ws.Range("AE3").Value = 1221
$ws.Range("AF3").Value = 0.88
$ws.Range("AG3").Value = 0
$ws.Range("AH3").Value = 0
$ws.Range("AI3").Value = 0
$ws.Range("AJ3").Value = 341155134
$ws.Range("J3").ClearContents()
$ws.Range("O3").ClearContents()
$ws.Range("U3").ClearContents()

# Row 4: 2016/12  (IFRS연결)
$ws.Range("D4").Value = 4671
$ws.Range("E4").Value = 79
$ws.Range("F4").Value = 79
$ws.Range("G4").Value = 168
$ws.Range("H4").Value = 114
$ws.Range("I4").Value = 114
$ws.Range("K4").Value = 39215
$ws.Range("L4").Value = 35064
$ws.Range("M4").Value = 4151
$ws.Range("N4").Value = 4151
$ws.Range("P4").Value = 1620
$ws.Range("Q4").Value = -1007
$ws.Range("R4").Value = -13
$ws.Range("S4").Value = 1559
$ws.Range("T4").Value = 10
$ws.Range("V4").Value = 6940
$ws.Range("W4").Value = 1.69
$ws.Range("X4").Value = 2.45
$ws.Range("Y4").Value = 2.77
$ws.Range("Z4").Value = 0.31
$ws.Range("AA4").Value = 844.76
$ws.Range("AB4").Value = 158.49
$ws.Range("AC4").Value = 33
$ws.Range("AD4").Value = 28.83
$ws.Range("AE4").Value = 1244
$ws.Range("AF4").Value = 0.77
$ws.Range("AG4").Value = 0
$ws.Range("AH4").Value = 0
$ws.Range("AI4").Value = 0
$ws.Range("AJ4").Value = 341155134
$ws.Range("J4").ClearContents()
$ws.Range("O4").ClearContents()
$ws.Range("U4").ClearContents()

# Row 5: 2017/12  (IFRS연결)
$ws.Range("D5").Value = 4673
$ws.Range("E5").Value = 285
$ws.Range("F5").Value = 285
$ws.Range("G5").Value = 254
$ws.Range("H5").Value = 188
$ws.Range("I5").Value = 188
$ws.Range("K5").Value = 37788
$ws.Range("L5").Value = 33418
$ws.Range("M5").Value = 4370
$ws.Range("N5").Value = 4370
$ws.Range("P5").Value = 1620
$ws.Range("Q5").Value = -81
$ws.Range("R5").Value = -367
$ws.Range("S5").Value = 220
$ws.Range("T5").Value = 129
$ws.Range("V5").Value = 8578
$ws.Range("W5").Value = 6.1
$ws.Range("X5").Value = 4.03
$ws.Range("Y5").Value = 4.42
$ws.Range("Z5").Value = 0.49
$ws.Range("AA5").Value = 764.75
$ws.Range("AB5").Value = 172.01
$ws.Range("AC5").Value = 55
$ws.Range("AD5").Value = 19.77
$ws.Range("AE5").Value = 1309
$ws.Range("AF5").Value = 0.82
$ws.Range("AG5").Value = 0
$ws.Range("AH5").Value = 0
$ws.Range("AI5").Value = 0
$ws.Range("AJ5").Value = 341155134
$ws.Range("J5").ClearContents()
$ws.Range("O5").ClearContents()
$ws.Range("U5").ClearContents()

# Row 6: 2018/12  (IFRS연결)
$ws.Range("D6").Value = 5343
$ws.Range("E6").Value = 129
$ws.Range("F6").Value = 129
$ws.Range("G6").Value = 213
$ws.Range("H6").Value = 139
$ws.Range("I6").Value = 139
$ws.Range("K6").Value = 37588
$ws.Range("L6").Value = 32173
$ws.Range("M6").Value = 5415
$ws.Range("N6").Value = 5415
$ws.Range("P6").Value = 2383
$ws.Range("Q6").Value = -2353
$ws.Range("R6").Value = 6
$ws.Range("S6").Value = 2493
$ws.Range("T6").Value = 45
$ws.Range("V6").Value = 7752
$ws.Range("W6").Value = 2.42
$ws.Range("X6").Value = 2.61
$ws.Range("Y6").Value = 2.85
$ws.Range("Z6").Value = 0.37
$ws.Range("AA6").Value = 594.12
$ws.Range("AB6").Value = 128.85
$ws.Range("AC6").Value = 39
$ws.Range("AD6").Value = 16.92
$ws.Range("AE6").Value = 1164
$ws.Range("AF6").Value = 0.57
$ws.Range("AG6").Value = 10
$ws.Range("AH6").Value = 1.52
$ws.Range("AI6").Value = 33.49
$ws.Range("AJ6").Value = 472590171
$ws.Range("U6").ClearContents()

# Rows 7-9: 2019E/2020E/2021E estimate columns removed; clear D:AJ, keep A-C labels.
$ws.Range("D7:AJ9").ClearContents()

